$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.961.55'
$ws.Range("E2").Value = '  +1.61%  '
$ws.Range("D3").Value = '3.330.27'
$ws.Range("E3").Value = '  +1.53%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '581.22'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.87%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.12'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.34%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("E8").Value = '  +1.42%  '
$ws.Range("D9").Value = '3.326.38'
$ws.Range("E9").Value = '  +1.56%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.183'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.581'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.81%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '46.96'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.61%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000273'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.14%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '686.60'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.05%  '
$ws.Range("D15").Value = '3.874.86'
$ws.Range("E15").Value = '  +1.77%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.43'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.13%  '
$ws.Range("D17").Value = '68.059.16'
$ws.Range("E17").Value = '  +1.62%  '
$ws.Range("E18").Value = '  -0.34%  '
$ws.Range("D19").Value = '3.330.39'
$ws.Range("E19").Value = '  +1.63%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.44'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.16%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.04'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.897'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.40%  '
$ws.Range("E23").Value = '  +4.94%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '17.04'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.30%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '99.08'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.91'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.27%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.69'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.47%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.58'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.67%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.04'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.93%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.57'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.83%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.07'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.09%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '566.54'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.15%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.01'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.99%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '57.47'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.81%  '
$ws.Range("E36").Value = '  -0.06%  '
$ws.Range("D37").Value = '3.706.80'
$ws.Range("E37").Value = '  -4.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.37'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.83%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '34.53'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +8.92%  '
$ws.Range("E40").Value = '  +3.39%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.19'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.89%  '
$ws.Range("E42").Value = '  +2.63%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.36'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.53%  '
$ws.Range("D44").Value = '0.0₃0675'
$ws.Range("E44").Value = '  +1.46%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.336'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.49%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0407'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.33%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.66'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +6.00%  '
$ws.Range("E48").Value = '  +1.37%  '
$ws.Range("E49").Value = '  -0.33%  '
$ws.Range("E50").Value = '  -2.80%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '130.06'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.05%  '
